$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.198.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'2.481.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'584.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'173.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.139"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").Value = "'4.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'2.933.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "'25.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'67.090.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "'2.395.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "'7.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "'350.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'4.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'1.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "'9.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'2.606.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'0.0₃0912"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "'505.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").Value = "'7.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'1.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'162.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "'18.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "'4.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'2.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("D45").Value = "'143.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "'0.0⁦0263"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  +0.41%  "
